$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "299.36"
Set-TextValue $ws.Range("E2") "-1.23%"
Set-TextValue $ws.Range("D3") "31.45"
Set-TextValue $ws.Range("E3") "-0.87%"
Set-TextValue $ws.Range("E4") "-1.44%"
Set-TextValue $ws.Range("D5") "0.07896"
Set-TextValue $ws.Range("E5") "1.01%"
Set-TextValue $ws.Range("D6") "2.265"
Set-TextValue $ws.Range("E6") "-6.83%"
Set-TextValue $ws.Range("D7") "7.811"
Set-TextValue $ws.Range("E7") "-1.74%"
Set-TextValue $ws.Range("D8") "3.851"
Set-TextValue $ws.Range("E8") "-0.50%"
Set-TextValue $ws.Range("D9") "0.9221"
Set-TextValue $ws.Range("E9") "1.29%"
Set-TextValue $ws.Range("D10") "0.1744"
Set-TextValue $ws.Range("E10") "0.91%"
Set-TextValue $ws.Range("D11") "0.07552"
Set-TextValue $ws.Range("E11") "3.01%"
Set-TextValue $ws.Range("D12") "0.09354"
Set-TextValue $ws.Range("E12") "14.91%"
Set-TextValue $ws.Range("D13") "0.03006"
Set-TextValue $ws.Range("E13") "-0.95%"
Set-TextValue $ws.Range("D14") "0.1002"
Set-TextValue $ws.Range("E14") "0.78%"
Set-TextValue $ws.Range("D15") "0.001505"
Set-TextValue $ws.Range("E15") "0.01%"
Set-TextValue $ws.Range("D16") "0.006070"
Set-TextValue $ws.Range("E16") "-0.38%"
Set-TextValue $ws.Range("D17") "3.475"
Set-TextValue $ws.Range("E17") "-0.66%"
Set-TextValue $ws.Range("D18") "2.243"
Set-TextValue $ws.Range("E18") "-0.04%"
Set-TextValue $ws.Range("E19") "0.87%"
Set-TextValue $ws.Range("D20") "0.1309"
Set-TextValue $ws.Range("E20") "-2.14%"
Set-TextValue $ws.Range("D21") "3.925"
Set-TextValue $ws.Range("E21") "-15.99%"
Set-TextValue $ws.Range("D22") "0.1711"
Set-TextValue $ws.Range("E22") "9.31%"
Set-TextValue $ws.Range("D23") "0.04615"
Set-TextValue $ws.Range("E23") "-0.81%"
Set-TextValue $ws.Range("D24") "0.001254"
Set-TextValue $ws.Range("E24") "-0.48%"
Set-TextValue $ws.Range("D25") "0.004475"
Set-TextValue $ws.Range("E25") "-1.00%"
Set-TextValue $ws.Range("D26") "0.0001248"
Set-TextValue $ws.Range("E26") "-7.48%"
Set-TextValue $ws.Range("D27") "0.0003397"
Set-TextValue $ws.Range("E27") "23.97%"
Set-TextValue $ws.Range("D39") "0.01733"
Set-TextValue $ws.Range("E39") "-2.82%"
Set-TextValue $ws.Range("D40") "0.04621"
Set-TextValue $ws.Range("E40") "1.46%"
Set-TextValue $ws.Range("D41") "0.006945"
Set-TextValue $ws.Range("E41") "-4.50%"
Set-TextValue $ws.Range("D42") "0.1358"
Set-TextValue $ws.Range("E42") "-0.34%"
Set-TextValue $ws.Range("D43") "0.002187"
Set-TextValue $ws.Range("E43") "-2.31%"
Set-TextValue $ws.Range("D44") "0.01030"
Set-TextValue $ws.Range("E44") "-4.23%"
Set-TextValue $ws.Range("D45") "0.00006274"
Set-TextValue $ws.Range("E45") "-3.12%"
Set-TextValue $ws.Range("E46") "0.05%"
Set-TextValue $ws.Range("D47") "0.007979"
Set-TextValue $ws.Range("E47") "-19.39%"
Set-TextValue $ws.Range("D48") "0.7451"
Set-TextValue $ws.Range("E48") "-9.20%"
Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("E49") "0.05%"
Set-TextValue $ws.Range("D50") "0.0002000"
Set-TextValue $ws.Range("E50") "0.05%"
